$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.557.07"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.082.68"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.56"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.87"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.438"
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.31"
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.383"
$ws.Range("E11").Value = "  +2.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.619.75"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.65"
$ws.Range("E14").Value = "  +3.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000166"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.605.11"
$ws.Range("E16").Value = "  +1.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.090.02"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.13"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.88"
$ws.Range("E19").Value = "  -1.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.10"
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "342.12"
$ws.Range("E21").Value = "  +1.49%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.504"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.75"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0914"
$ws.Range("E27").Value = "  -1.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.58"
$ws.Range("E28").Value = "  +2.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.22"
$ws.Range("E29").Value = "  +1.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.84"
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.99"
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("E32").Value = "  +2.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.29"
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.59"
$ws.Range("E34").Value = "  +1.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.09"
$ws.Range("E35").Value = "  +3.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.88"
$ws.Range("E36").Value = "  -3.65%  "
$ws.Range("E37").Value = "  +4.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0679"
$ws.Range("E38").Value = "  -1.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.131.56"
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.89"
$ws.Range("E40").Value = "  +0.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.77"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("E42").Value = "  +7.30%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("E44").Value = "  -0.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.270.87"
$ws.Range("E45").Value = "  -0.65%  "
$ws.Range("E46").Value = "  +1.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.70"
$ws.Range("E47").Value = "  +1.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.957"
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.99"
$ws.Range("E49").Value = "  +1.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.744"
$ws.Range("E50").Value = "  +7.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "263.57"
$ws.Range("E51").Value = "  +10.22%  "
